$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("gggu", "gjhh", "a", "a", "a", "a", "a", "a", "a", "a", "a", "a", 50,  "Reprovado(a)"),
    @("gggu", "gjhh", "a", "b", "d", "c", "c", "b", "d", "a", "b", "a", 110, "Aprovado(a)"),
    @("gggu", "gjhh", "b", "b", "a", "a", "c", "b", "a", "a", "b", "a", 190, "Aprovado(a)"),
    @("ff",   "uuu",  "a", "b", "c", "a", "a", "c", "b", "a", "b", "b", 60,  "Aprovado(a)")
)

$startRow = 17
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
